$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3-Year Summary")

# Row 5: Software Licenses -> Support & Maintenance
$ws.Range("A5").Value = "Support & Maintenance"

# Row 6: Support & Maintenance -> Risk Mitigation
$ws.Range("A6").Value = "Risk Mitigation"

# Row 7: TOTAL -> Operational Savings, with new SUMIF formulas
$ws.Range("A7").Value = "Operational Savings"
$ws.Range("B7").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A7,'Infrastructure Costs'!`$G:`$G)"
$ws.Range("C7").Formula = "=SUMIF(Credits!`$A:`$A,A7,Credits!`$C:`$C)"
$ws.Range("D7").Formula = "=B7+C7"
$ws.Range("E7").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A7,'Infrastructure Costs'!`$H:`$H)"
$ws.Range("F7").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A7,'Infrastructure Costs'!`$I:`$I)"
$ws.Range("G7").Formula = "=D7+E7+F7"

# Row 8: Password Reset Savings -> Compliance Enablement
$ws.Range("A8").Value = "Compliance Enablement"

# Row 9: Net Investment After Savings -> TOTAL, with new SUM formulas
$ws.Range("A9").Value = "TOTAL"
$ws.Range("B9").Formula = "=SUM(B3:B8)"
$ws.Range("C9").Formula = "=SUM(C3:C8)"
$ws.Range("D9").Formula = "=SUM(D3:D8)"
$ws.Range("E9").Formula = "=SUM(E3:E8)"
$ws.Range("F9").Formula = "=SUM(F3:F8)"
$ws.Range("G9").Formula = "=SUM(G3:G8)"

# The old data below the former TOTAL row is gone now that TOTAL moved to row 9;
# register row 10 as an empty trailing row (clears any stray formatting/content there).
$ws.Rows.Item(10).OutlineLevel = 1
$ws.Rows.Item(10).OutlineLevel = 0
